$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 30, shifting existing rows 30-76 down to 31-77
$ws.Rows(30).Insert()

# Fill in the new row 30 with the new data record
$ws.Range("A30").Value = 5
$ws.Range("B30").Value = "Macroferia Regional de Talca"
$ws.Range("C30").Value = "Maule"
$ws.Range("D30").Value = 44469
$ws.Range("D30").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E30").Value = 7
$ws.Range("F30").Value = "Fruta"
$ws.Range("G30").Value = 100108
$ws.Range("H30").Value = "Tropicales y subtropicales"
$ws.Range("I30").Value = 100108002
$ws.Range("J30").Value = "Mango"
$ws.Range("K30").Value = "Sin especificar"
$ws.Range("L30").Value = "Primera"
$ws.Range("M30").Value = 108
$ws.Range("N30").Value = 8000
$ws.Range("O30").Value = 8000
$ws.Range("P30").Value = 8000
$ws.Range("Q30").Value = "$/bandeja 4 kilos"
$ws.Range("R30").Value = "Brasil"
$ws.Range("S30").Value = 2000
$ws.Range("T30").Value = 4
